$wb = $excel.ActiveWorkbook

# ALC row 6
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(6, 8).Value = 211.66667  # H6: was 83.5
$ws.Cells.Item(6, 9).Value = 267.5  # I6: was 81.14286
$ws.Cells.Item(6, 11).Value = 802.5  # K6: was 243.42858
$ws.Cells.Item(6, 13).Value = -690.5  # M6: was -131.42858

# ALC row 33
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 361.14285  # H33: was 342.06668
$ws.Cells.Item(33, 10).Value = 508.375  # J33: was 421.7
$ws.Cells.Item(33, 12).Value = 508.375  # L33: was 421.7
$ws.Cells.Item(33, 14).Value = -966.375  # N33: was -879.7

# ALC row 86
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(86, 8).Value = 2501085.5  # H86: was 3334299.2
$ws.Cells.Item(86, 10).Value = 1444  # J86: was 0
$ws.Cells.Item(86, 12).Value = 1444  # L86: was 0
$ws.Cells.Item(86, 14).Value = -3690  # N86: was None

# ALC row 89
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(89, 8).Value = 2501085.5  # H89: was 3334299.2
$ws.Cells.Item(89, 10).Value = 1444  # J89: was 0
$ws.Cells.Item(89, 12).Value = 7220  # L89: was 0
$ws.Cells.Item(89, 14).Value = -18452  # N89: was None

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value = 4747.722  # H137: was 4903.706
$ws.Cells.Item(137, 9).Value = 3403.96  # I137: was 3379.2083
$ws.Cells.Item(137, 10).Value = 7801.727  # J137: was 8562.5
$ws.Cells.Item(137, 11).Value = 10211.88  # K137: was 10137.6249
$ws.Cells.Item(137, 12).Value = 23405.181  # L137: was 25687.5
$ws.Cells.Item(137, 13).Value = -7661.880000000001  # M137: was -7587.624899999999
$ws.Cells.Item(137, 14).Value = -28505.181  # N137: was -30787.5

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(138, 8).Value = 3409.4385  # H138: was 3391.2932
$ws.Cells.Item(138, 10).Value = 4824.6  # J138: was 4756.0557
$ws.Cells.Item(138, 12).Value = 14473.8  # L138: was 14268.1671
$ws.Cells.Item(138, 14).Value = -24753.8  # N138: was -24548.1671

# ARM row 17
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(17, 8).Value = 1100  # H17: was 0
$ws.Cells.Item(17, 9).Value = 1100  # I17: was 0
$ws.Cells.Item(17, 11).Value = 1100  # K17: was 0
$ws.Cells.Item(17, 13).Value = -927  # M17: was None

# ARM row 25
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(25, 8).Value = 7228.75  # H25: was 8753
$ws.Cells.Item(25, 9).Value = 2971.6667  # I25: was 3015
$ws.Cells.Item(25, 10).Value = 20000  # J25: was 10665.667
$ws.Cells.Item(25, 11).Value = 2971.6667  # K25: was 3015
$ws.Cells.Item(25, 12).Value = 20000  # L25: was 10665.667
$ws.Cells.Item(25, 13).Value = -2569.6667  # M25: was -2613
$ws.Cells.Item(25, 14).Value = -20804  # N25: was -11469.667

# ARM row 35
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(35, 8).Value = 10754.4  # H35: was 10085.444
$ws.Cells.Item(35, 9).Value = 10754.4  # I35: was 8846.125
$ws.Cells.Item(35, 10).Value = 0  # J35: was 20000
$ws.Cells.Item(35, 11).Value = 10754.4  # K35: was 8846.125
$ws.Cells.Item(35, 12).Value = 0  # L35: was 20000
$ws.Cells.Item(35, 13).ClearContents()  # M35: was -8440.125
$ws.Cells.Item(35, 14).Value = -10348.4  # N35: was -20812

# ARM row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 4486.075  # H61: was 4576.757
$ws.Cells.Item(61, 9).Value = 4518.8335  # I61: was 4658.5
$ws.Cells.Item(61, 10).Value = 4459.273  # J61: was 4514.476
$ws.Cells.Item(61, 11).Value = 4518.8335  # K61: was 4658.5
$ws.Cells.Item(61, 12).Value = 4459.273  # L61: was 4514.476
$ws.Cells.Item(61, 13).Value = -4306.8335  # M61: was -4446.5
$ws.Cells.Item(61, 14).Value = -4883.273  # N61: was -4938.476

# ARM row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(136, 8).Value = 4486.075  # H136: was 4576.757
$ws.Cells.Item(136, 9).Value = 4518.8335  # I136: was 4658.5
$ws.Cells.Item(136, 10).Value = 4459.273  # J136: was 4514.476
$ws.Cells.Item(136, 11).Value = 13556.5005  # K136: was 13975.5
$ws.Cells.Item(136, 12).Value = 13377.819  # L136: was 13543.428
$ws.Cells.Item(136, 13).Value = -11006.5005  # M136: was -11425.5
$ws.Cells.Item(136, 14).Value = -18477.819  # N136: was -18643.428

# BSM row 25
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(25, 8).Value = 7514  # H25: was 4213.5
$ws.Cells.Item(25, 9).Value = 7514  # I25: was 4213.5
$ws.Cells.Item(25, 11).Value = 7514  # K25: was 4213.5
$ws.Cells.Item(25, 13).Value = -7279  # M25: was -3978.5

# BSM row 37
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(37, 8).Value = 5303.8335  # H37: was 5570.3335
$ws.Cells.Item(37, 9).Value = 1364.6  # I37: was 1684.4
$ws.Cells.Item(37, 11).Value = 1364.6  # K37: was 1684.4
$ws.Cells.Item(37, 13).Value = -1227.6  # M37: was -1547.4

# BSM row 81
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(81, 8).Value = 10058.333  # H81: was 10876
$ws.Cells.Item(81, 10).Value = 10058.333  # J81: was 10876
$ws.Cells.Item(81, 12).Value = 10058.333  # L81: was 10876
$ws.Cells.Item(81, 14).Value = -12180.333  # N81: was -12998

# BSM row 84
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(84, 8).Value = 10058.333  # H84: was 10876
$ws.Cells.Item(84, 10).Value = 10058.333  # J84: was 10876
$ws.Cells.Item(84, 12).Value = 30174.999  # L84: was 32628
$ws.Cells.Item(84, 14).Value = -40782.999  # N84: was -43236

# BSM row 95
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(95, 8).Value = 11955  # H95: was 12121
$ws.Cells.Item(95, 10).Value = 11955  # J95: was 12121
$ws.Cells.Item(95, 12).Value = 11955  # L95: was 12121
$ws.Cells.Item(95, 14).Value = -17447  # N95: was -17613

# BSM row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 11166723  # H107: was 4785742
$ws.Cells.Item(107, 9).Value = 12561626  # I107: was 5024954.5
$ws.Cells.Item(107, 10).Value = 7500  # J107: was 1500
$ws.Cells.Item(107, 11).Value = 12561626  # K107: was 5024954.5
$ws.Cells.Item(107, 12).Value = 7500  # L107: was 1500
$ws.Cells.Item(107, 13).Value = -12559706  # M107: was -5023034.5
$ws.Cells.Item(107, 14).Value = -11340  # N107: was -5340

# CRP row 12
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(12, 8).Value = 0  # H12: was 300
$ws.Cells.Item(12, 10).Value = 0  # J12: was 300
$ws.Cells.Item(12, 12).ClearContents()  # L12: was 300
$ws.Cells.Item(12, 14).Value = 0  # N12: was -640

# CRP row 51
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(51, 8).Value = 175000  # H51: was 122500
$ws.Cells.Item(51, 10).Value = 0  # J51: was 70000
$ws.Cells.Item(51, 12).Value = 0  # L51: was 70000
$ws.Cells.Item(51, 14).ClearContents()  # N51: was -71472

# CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 66671810  # H58: was 55560260
$ws.Cells.Item(58, 9).Value = 100003480  # I58: was 83336670
$ws.Cells.Item(58, 10).Value = 8461  # J58: was 7432.5
$ws.Cells.Item(58, 11).Value = 100003480  # K58: was 83336670
$ws.Cells.Item(58, 12).Value = 8461  # L58: was 7432.5
$ws.Cells.Item(58, 13).Value = -100003277  # M58: was -83336467
$ws.Cells.Item(58, 14).Value = -8867  # N58: was -7838.5

# CRP row 61
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(61, 8).Value = 175000  # H61: was 122500
$ws.Cells.Item(61, 10).Value = 0  # J61: was 70000
$ws.Cells.Item(61, 12).Value = 0  # L61: was 70000
$ws.Cells.Item(61, 14).ClearContents()  # N61: was -70696

# CRP row 62
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(62, 8).Value = 7735  # H62: was 11775.667
$ws.Cells.Item(62, 9).Value = 1900.25  # I62: was 1904
$ws.Cells.Item(62, 10).Value = 11624.833  # J62: was 13750
$ws.Cells.Item(62, 11).Value = 1900.25  # K62: was 1904
$ws.Cells.Item(62, 12).Value = 11624.833  # L62: was 13750
$ws.Cells.Item(62, 13).Value = -1276.25  # M62: was -1280
$ws.Cells.Item(62, 14).Value = -12872.833  # N62: was -14998

# CRP row 65
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(65, 8).Value = 7735  # H65: was 11775.667
$ws.Cells.Item(65, 9).Value = 1900.25  # I65: was 1904
$ws.Cells.Item(65, 10).Value = 11624.833  # J65: was 13750
$ws.Cells.Item(65, 11).Value = 9501.25  # K65: was 9520
$ws.Cells.Item(65, 12).Value = 58124.165  # L65: was 68750
$ws.Cells.Item(65, 13).Value = -6381.25  # M65: was -6400
$ws.Cells.Item(65, 14).Value = -64364.165  # N65: was -74990

# CRP row 88
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(88, 8).Value = 39933  # H88: was 49950
$ws.Cells.Item(88, 10).Value = 39933  # J88: was 49950
$ws.Cells.Item(88, 12).Value = 39933  # L88: was 49950
$ws.Cells.Item(88, 14).Value = -40745  # N88: was -50762

# CRP row 91
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(91, 8).Value = 39933  # H91: was 49950
$ws.Cells.Item(91, 10).Value = 39933  # J91: was 49950
$ws.Cells.Item(91, 12).Value = 39933  # L91: was 49950
$ws.Cells.Item(91, 14).Value = -42741  # N91: was -52758

# CRP row 95
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(95, 8).Value = 21849.75  # H95: was 13049.857
$ws.Cells.Item(95, 10).Value = 21849.75  # J95: was 13049.857
$ws.Cells.Item(95, 12).Value = 21849.75  # L95: was 13049.857
$ws.Cells.Item(95, 14).Value = -27341.75  # N95: was -18541.857

# CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(136, 8).Value = 66671810  # H136: was 55560260
$ws.Cells.Item(136, 9).Value = 100003480  # I136: was 83336670
$ws.Cells.Item(136, 10).Value = 8461  # J136: was 7432.5
$ws.Cells.Item(136, 11).Value = 300010440  # K136: was 250010010
$ws.Cells.Item(136, 12).Value = 25383  # L136: was 22297.5
$ws.Cells.Item(136, 13).Value = -300007890  # M136: was -250007460
$ws.Cells.Item(136, 14).Value = -30483  # N136: was -27397.5

# CUL row 17
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(17, 8).Value = 322.75  # H17: was 469.25
$ws.Cells.Item(17, 9).Value = 322.75  # I17: was 469.25
$ws.Cells.Item(17, 11).Value = 968.25  # K17: was 1407.75
$ws.Cells.Item(17, 13).Value = -799.25  # M17: was -1238.75

# CUL row 21
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(21, 8).Value = 3999  # H21: was 1884
$ws.Cells.Item(21, 9).Value = 0  # I21: was 768.5
$ws.Cells.Item(21, 10).Value = 3999  # J21: was 2999.5
$ws.Cells.Item(21, 11).Value = 0  # K21: was 2305.5
$ws.Cells.Item(21, 12).ClearContents()  # L21: was 8998.5
$ws.Cells.Item(21, 13).Value = 11997  # M21: was -2132.5
$ws.Cells.Item(21, 14).Value = -12343  # N21: was -9344.5

# CUL row 37
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(37, 8).Value = 137495.5  # H37: was 134368.88
$ws.Cells.Item(37, 10).Value = 137495.5  # J37: was 134368.88
$ws.Cells.Item(37, 12).Value = 412486.5  # L37: was 403106.64
$ws.Cells.Item(37, 14).Value = -412710.5  # N37: was -403330.64

# CUL row 140
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(140, 8).Value = 26043476  # H140: was 24511628
$ws.Cells.Item(140, 9).Value = 30865782  # I140: was 28737250
$ws.Cells.Item(140, 11).Value = 92597346  # K140: was 86211750
$ws.Cells.Item(140, 13).Value = -92592166  # M140: was -86206570

# GSM row 14
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(14, 8).Value = 5565031.5  # H14: was 8743821
$ws.Cells.Item(14, 9).Value = 5565031.5  # I14: was 10200292
$ws.Cells.Item(14, 10).Value = 0  # J14: was 5000
$ws.Cells.Item(14, 11).Value = 5565031.5  # K14: was 10200292
$ws.Cells.Item(14, 12).Value = 0  # L14: was 5000
$ws.Cells.Item(14, 13).ClearContents()  # M14: was -10200124
$ws.Cells.Item(14, 14).Value = -5564863.5  # N14: was -5336

# GSM row 130
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(130, 8).Value = 0  # H130: was 48000
$ws.Cells.Item(130, 10).Value = 0  # J130: was 48000
$ws.Cells.Item(130, 12).ClearContents()  # L130: was 48000
$ws.Cells.Item(130, 14).Value = 0  # N130: was -58040

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 142860860  # H132: was 166670200
$ws.Cells.Item(132, 9).Value = 200002200  # I132: was 200002750
$ws.Cells.Item(132, 10).Value = 7499.5  # J132: was 7500
$ws.Cells.Item(132, 11).Value = 600006600  # K132: was 600008250
$ws.Cells.Item(132, 12).Value = 22498.5  # L132: was 22500
$ws.Cells.Item(132, 13).Value = -600004070  # M132: was -600005720
$ws.Cells.Item(132, 14).Value = -27558.5  # N132: was -27560

# LTW row 16
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 5987.7393  # H16: was 6250.7393
$ws.Cells.Item(16, 9).Value = 4773.4  # I16: was 5188.4
$ws.Cells.Item(16, 10).Value = 14083.333  # J16: was 13333
$ws.Cells.Item(16, 11).Value = 4773.4  # K16: was 5188.4
$ws.Cells.Item(16, 12).Value = 14083.333  # L16: was 13333
$ws.Cells.Item(16, 13).Value = -4603.4  # M16: was -5018.4
$ws.Cells.Item(16, 14).Value = -14423.333  # N16: was -13673

# LTW row 17
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(17, 8).Value = 17001666  # H17: was 392873180
$ws.Cells.Item(17, 10).Value = 1000000  # J17: was 497979070
$ws.Cells.Item(17, 12).Value = 1000000  # L17: was 497979070
$ws.Cells.Item(17, 14).Value = -1000340  # N17: was -497979410

# LTW row 140
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(140, 8).Value = 78999.5  # H140: was 30429
$ws.Cells.Item(140, 10).Value = 78999.5  # J140: was 30429
$ws.Cells.Item(140, 12).Value = 78999.5  # L140: was 30429
$ws.Cells.Item(140, 14).Value = -89359.5  # N140: was -40789

# WVR row 41
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(41, 8).Value = 9222.5  # H41: was 9269
$ws.Cells.Item(41, 9).Value = 9176  # I41: was 0
$ws.Cells.Item(41, 11).Value = 9176  # K41: was 0
$ws.Cells.Item(41, 13).Value = -8786  # M41: was None

# WVR row 62
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 26166.334  # H62: was 18299.4
$ws.Cells.Item(62, 9).Value = 32750  # I62: was 14916
$ws.Cells.Item(62, 10).Value = 22874.5  # J62: was 23374.5
$ws.Cells.Item(62, 11).Value = 32750  # K62: was 14916
$ws.Cells.Item(62, 12).Value = 22874.5  # L62: was 23374.5
$ws.Cells.Item(62, 13).Value = -32126  # M62: was -14292
$ws.Cells.Item(62, 14).Value = -24122.5  # N62: was -24622.5

# WVR row 65
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(65, 8).Value = 26166.334  # H65: was 18299.4
$ws.Cells.Item(65, 9).Value = 32750  # I65: was 14916
$ws.Cells.Item(65, 10).Value = 22874.5  # J65: was 23374.5
$ws.Cells.Item(65, 11).Value = 163750  # K65: was 74580
$ws.Cells.Item(65, 12).Value = 114372.5  # L65: was 116872.5
$ws.Cells.Item(65, 13).Value = -160630  # M65: was -71460
$ws.Cells.Item(65, 14).Value = -120612.5  # N65: was -123112.5

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 9369.5  # H132: was 10294.154
$ws.Cells.Item(132, 9).Value = 8190.5  # I132: was 8256.333000000001
$ws.Cells.Item(132, 10).Value = 17622.5  # J132: was 34748
$ws.Cells.Item(132, 11).Value = 24571.5  # K132: was 24768.999
$ws.Cells.Item(132, 12).Value = 52867.5  # L132: was 104244
$ws.Cells.Item(132, 13).Value = -22041.5  # M132: was -22238.999
$ws.Cells.Item(132, 14).Value = -57927.5  # N132: was -109304
